$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

$ws.Range("D5").Value = "KONTOSTAND AM 20.06.2024"

$ws.Range("B6").Value = "21.06."
$ws.Range("C6").Value = "22.06."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 65244081"
$ws.Range("E6").Value = "85,43-"

$ws.Range("B7").Value = "25.06."
$ws.Range("C7").Value = "26.06."
$ws.Range("D7").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E7").Value = "65,61-"

$ws.Range("B8").Value = "26.06."
$ws.Range("C8").Value = "27.06."
$ws.Range("D8").Value = "BEITRAG Allianz SE K-24405313"
$ws.Range("E8").Value = "54,88-"

$ws.Range("D12").Value = "KONTOSTAND AM 30.06.2024"
$ws.Range("E12").Value = "205,92-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 07.07.2024"
